$d = $word.ActiveDocument

# Locate the end of the sentence "... у нашу папку. " to insert the new runs right after it.
$rng = $d.Content
$found = $rng.Find.Execute("у нашу папку. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Collapse the found range to its end point, so subsequent inserts land right after it.
$rng.Collapse(0)

# First new run: Ukrainian text, sz 28.
$rng.InsertAfter("У ході виконання лабораторної роботи, виникала проблема конфліктів. Поки я виконувала її, іншими одногрупниками були додані файли, тому потрібно було використовувати операцію ")
$rng.Font.Size = 14
$rng.LanguageID = 1058

# Move to the end of what was just inserted.
$rng.Collapse(0)

# Second new run: English "merge", sz 28.
$rng.InsertAfter("merge")
$rng.Font.Size = 14
$rng.LanguageID = 1033

$rng.Collapse(0)

# Third new run: Ukrainian ".", sz 28.
$rng.InsertAfter(".")
$rng.Font.Size = 14
$rng.LanguageID = 1058
